# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as plain text so values such as
# "4.272" or "19.60" are not silently re-interpreted as numbers
# (they are formatted price strings, some using "." as a thousands
# separator, e.g. "26.069.14").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.069.14"
$ws.Range("D3").Value = "1.648.02"
$ws.Range("D5").Value = "215.33"
$ws.Range("D6").Value = "0.5095"
$ws.Range("D8").Value = "0.2586"
$ws.Range("D9").Value = "0.06423"
$ws.Range("D10").Value = "19.60"
$ws.Range("D11").Value = "0.07721"
$ws.Range("D12").Value = "1.648.24"
$ws.Range("D13").Value = "4.272"
$ws.Range("D14").Value = "1.877.47"
$ws.Range("D15").Value = "0.5461"
$ws.Range("D16").Value = "0.0₅7982"
$ws.Range("D17").Value = "63.80"
$ws.Range("D18").Value = "26.099.12"
$ws.Range("D20").Value = "207.06"
$ws.Range("D21").Value = "4.379"
$ws.Range("D22").Value = "10.03"
$ws.Range("D23").Value = "5.999"
$ws.Range("D25").Value = "1.869"
$ws.Range("D26").Value = "143.13"
$ws.Range("D27").Value = "0.1166"
$ws.Range("D28").Value = "6.902"
$ws.Range("D29").Value = "15.78"
$ws.Range("D30").Value = "0.05065"
$ws.Range("D31").Value = "1.243"
$ws.Range("D32").Value = "3.329"
$ws.Range("D33").Value = "3.225"
$ws.Range("D34").Value = "1.547"
$ws.Range("D35").Value = "2.347"
$ws.Range("D36").Value = "0.9131"
$ws.Range("D37").Value = "2.643"
$ws.Range("D38").Value = "0.5700"
$ws.Range("D39").Value = "1.143.43"
$ws.Range("D40").Value = "0.01572"
$ws.Range("D43").Value = "5.659"
$ws.Range("D44").Value = "0.8229"
$ws.Range("D45").Value = "99.80"
$ws.Range("D46").Value = "1.789.06"
$ws.Range("D47").Value = "0.0₈113"
$ws.Range("D48").Value = "0.4538"
$ws.Range("D49").Value = "1.010"
$ws.Range("D50").Value = "55.15"
$ws.Range("D51").Value = "7.801"

$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("E10").Value = "  -4.68%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("E17").Value = "  -4.70%  "
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("E21").Value = "  -4.68%  "
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +8.14%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  -3.54%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("E51").Value = "  -3.14%  "

